# setProperty.xlsx update: add a new column (setProperty/json/age/name split),
# migrate JSON sample cells, update rich-text coloring, column widths & selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Colors (BGR decimal equivalents of the target RGB hex) ---
$colBlack  = 0        # 000000 - default/base text color (was "indexed 8")
$colDark   = 526344   # 080808 - punctuation / plain text
$colPurple = 9703559  # 871094 - JSON keys
$colGreen  = 1539334  # 067D17 - JSON values

$fontName = "游ゴシック"

# --- 1. Insert a new column before C, shifting old C (echo/${file}) -> D
#        and old D (wait/2000) -> E. Formatting of the new column is copied
#        from the column to its left (B), matching the target style pattern. ---
$ws.Columns.Item(3).Insert()

# --- 2. Row 1 (headers) ---
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "setProperty"
$ws.Range("C1").Value = "setProperty"
$ws.Range("D1").Value = "echo"
$ws.Range("E1").Value = "wait"

# --- 3. Row 2 ---
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = '{"type":"json"}'
$ws.Range("C2").Value = "json"
$ws.Range("D2").Value = '${age}_${name}'
$ws.Range("E2").Value = ""

# --- 4. Row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = '{"age":"40"}'
$ws.Range("C3").Value = '{"name":"hugang"}'
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = 2000

# --- 5. Base fonts per cell (plain, non-rich-text cells first) ---
foreach ($addr in @("A1","B1","C1","D1","A2","A3")) {
  $c = $ws.Range($addr)
  $c.Font.Name = $fontName
  $c.Font.Color = $colBlack
}
foreach ($addr in @("E1","E3")) {
  $c = $ws.Range($addr)
  $c.Font.Name = $fontName
  $c.Font.Color = $colBlack
}
$e2 = $ws.Range("E2")
$e2.Font.Name = $fontName
$e2.Font.Color = $colBlack

foreach ($addr in @("C2","D2")) {
  $c = $ws.Range($addr)
  $c.Font.Name = $fontName
  $c.Font.Color = $colDark
}

# --- 6. Rich-text JSON cells: B2, B3, C3 ---
# B2 = {"type":"json"}
$b2 = $ws.Range("B2")
$b2.Font.Name = $fontName
$b2.Font.Color = $colDark
$b2.Characters(2,6).Font.Color = $colPurple
$b2.Characters(2,6).Font.Name = $fontName
$b2.Characters(8,1).Font.Color = $colDark
$b2.Characters(8,1).Font.Name = $fontName
$b2.Characters(9,6).Font.Color = $colGreen
$b2.Characters(9,6).Font.Name = $fontName
$b2.Characters(15,1).Font.Color = $colDark
$b2.Characters(15,1).Font.Name = $fontName

# B3 = {"age":"40"}
$b3 = $ws.Range("B3")
$b3.Font.Name = $fontName
$b3.Font.Color = $colDark
$b3.Characters(2,5).Font.Color = $colPurple
$b3.Characters(2,5).Font.Name = $fontName
$b3.Characters(7,1).Font.Color = $colDark
$b3.Characters(7,1).Font.Name = $fontName
$b3.Characters(8,4).Font.Color = $colGreen
$b3.Characters(8,4).Font.Name = $fontName
$b3.Characters(12,1).Font.Color = $colDark
$b3.Characters(12,1).Font.Name = $fontName

# C3 = {"name":"hugang"}
$c3 = $ws.Range("C3")
$c3.Font.Name = $fontName
$c3.Font.Color = $colDark
$c3.Characters(2,6).Font.Color = $colPurple
$c3.Characters(2,6).Font.Name = $fontName
$c3.Characters(8,1).Font.Color = $colDark
$c3.Characters(8,1).Font.Name = $fontName
$c3.Characters(9,8).Font.Color = $colGreen
$c3.Characters(9,8).Font.Name = $fontName
$c3.Characters(17,1).Font.Color = $colDark
$c3.Characters(17,1).Font.Name = $fontName

# --- 7. Row heights: re-fit row 1 back down to the default (content is
#        short now, no more wrapping needed). ---
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# --- 8. Column widths (best-fit, approximated through the runtime's pixel
#        quantized ColumnWidth setter). ---
$ws.Columns.Item(1).ColumnWidth = 10.65
$ws.Columns.Item(2).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(4).ColumnWidth = 13.79
$ws.Columns.Item(5).ColumnWidth = 5.5

# --- 9. Selection state ---
[void]$ws.Range("B4").Select()
